# Update "Ambati Rayudu" Chennai Super Kings batting-innings log.
# The sheet's existing rows were reshuffled and one new innings (row 12)
# was appended, per the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure numeric-looking values in columns C:F stay stored as text,
# matching the source data (runs/balls/fours/sixes were text-typed cells).
$ws.Range("C2:F12").NumberFormat = "@"

# Final target grid for A2:F12 (playerName, teamName, runs, balls, fours, sixes)
$data = @(
    @("Ambati Rayudu ", "Chennai Super Kings", "2",  "3",  "0", "0"),
    @("Ambati Rayudu ", "Chennai Super Kings", "38", "20", "5", "1"),
    @("Ambati Rayudu ", "Chennai Super Kings", "30", "30", "2", "0"),
    @("Ambati Rayudu ", "Chennai Super Kings", "39", "27", "3", "2"),
    @("Ambati Rayudu ", "Chennai Super Kings", "30", "27", "3", "0"),
    @("Ambati Rayudu ", "Chennai Super Kings", "13", "19", "2", "0"),
    @("Ambati Rayudu ", "Chennai Super Kings", "41", "34", "3", "2"),
    @("Ambati Rayudu ", "Chennai Super Kings", "45", "25", "1", "4"),
    @("Ambati Rayudu ", "Chennai Super Kings", "8",  "9",  "1", "0"),
    @("Ambati Rayudu ", "Chennai Super Kings", "71", "48", "6", "3"),
    @("Ambati Rayudu ", "Chennai Super Kings", "42", "40", "4", "0")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $row = $row + 1
}
